# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Row 2 (previously LINE_B) now reports the LINE_A figures, and the former
# LINE_B figures move to a newly-appended row 3 with their own recomputed
# quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# --- Update existing row 2: MAT_B/LINE_B -> MAT_A/LINE_A, revised quantities ---
$ws.Range("A2").Value = "MAT_A"
$ws.Range("C2").Value = "LINE_A"
$ws.Range("G2").Value = 870
$ws.Range("H2").Value = 870
$ws.Range("J2").Value = 827

# --- Append new row 3: the original MAT_B / LINE_B line, recomputed ---
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"

$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("D2:F2").Copy()
$ws.Range("D3:F3").PasteSpecial(-4122)

$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 120
$ws.Range("I2").Copy($ws.Range("I3"))
$ws.Range("J3").Value = 106
